# DPLKINV021-001 - Setup Custodian Investasi - General Tambah Data
# Update the "Kode Custody" / "Nama Custody" test data from QAS / QAS TEST
# to QAST / QAS TESTING, and refresh the summary cell (F2) that
# concatenates all the field values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
    "Username : 31246;",
    "Password : bni1234;",
    "Role : 20/21 - Analis Investasi;",
    "Kode Custody : QAST;",
    "Nama Custody : QAS TESTING;",
    "Alamat1 : Jalan Testing;",
    "Alamat2 : -;",
    "Kota : JAKARTA;",
    "Telp : 99999999;",
    "Email : test@gmail.com;",
    "PIC : -;",
    "PIC2 : -;",
    "Kode BI : 009 : Bank Negara Indonesia 1946;",
    "Nama Bank : Bank Negara Indonesia 1946;",
    "Cabang Bank : Pejompongan;",
    "No. Rekening : 12131313131;",
    "Pemilik Rekening : Tester;",
    "No. Giro : -;",
    "Status Register : 1 : Lanjutkan ke Verifikasi;",
    "Keterangan Register : DATA UNTUK DIVERIFIKASI"
)
$summary = [string]::Join("`n", $lines)

# Set M2 (KODE_CUSTODY) before F2 so newly-minted shared strings land in
# the same order as the authored workbook (QAST, then the summary text,
# then the NAMA_CUSTODY value).
$ws.Range("M2").Value = "QAST"
$ws.Range("F2").Value = $summary
$ws.Range("N2").Value = "QAS TESTING"

# Matches the author's final selection in the saved file.
$ws.Range("O2").Select()
